$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1252.9642
$ws.Range("J17").Value = 1394.9131
$ws.Range("L17").Value = 4184.7393
$ws.Range("N17").Value = -4520.7393
$ws.Range("H38").Value = 569
$ws.Range("I38").Value = 103.5
$ws.Range("J38").Value = 1500
$ws.Range("K38").Value = 310.5
$ws.Range("L38").Value = 4500
$ws.Range("M38").Value = 61.5
$ws.Range("N38").Value = -5244
$ws.Range("H40").Value = 4000
$ws.Range("I40").Value = 4000
$ws.Range("K40").Value = 4000
$ws.Range("M40").Value = -3825
$ws.Range("H64").Value = 21665.666
$ws.Range("I64").Value = 21665.666
$ws.Range("K64").Value = 21665.666
$ws.Range("M64").Value = -21417.666
$ws.Range("H67").Value = 21665.666
$ws.Range("I67").Value = 21665.666
$ws.Range("K67").Value = 21665.666
$ws.Range("M67").Value = -20807.666
$ws.Range("H116").Value = 11002.5
$ws.Range("I116").Value = 11002.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 11002.5
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -7560.5
$ws.Range("N116").ClearContents()
$ws.Range("H137").Value = 2740
$ws.Range("I137").Value = 2237.625
$ws.Range("J137").Value = 3105.3635
$ws.Range("K137").Value = 6712.875
$ws.Range("L137").Value = 9316.0905
$ws.Range("M137").Value = -4162.875
$ws.Range("N137").Value = -14416.0905

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3988.889
$ws.Range("I132").Value = 3782.4
$ws.Range("J132").Value = 4247
$ws.Range("K132").Value = 11347.2
$ws.Range("L132").Value = 12741
$ws.Range("M132").Value = -8817.200000000001
$ws.Range("N132").Value = -17801

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 20399.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3908.25
$ws.Range("I31").Value = 3908.25
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 3908.25
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -3613.25
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 3908.25
$ws.Range("I34").Value = 3908.25
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 3908.25
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -3706.25
$ws.Range("N34").ClearContents()
$ws.Range("H56").Value = 10000
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 10000
$ws.Range("M56").Value = -9155
$ws.Range("H62").Value = 16000
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 30000
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 30000
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -31248
$ws.Range("H65").Value = 16000
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 30000
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 150000
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -156240

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 1174.8334
$ws.Range("I134").Value = 1174.8334
$ws.Range("K134").Value = 3524.5002
$ws.Range("M134").Value = 1545.4998

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 474021.8
$ws.Range("I2").Value = 150737.64
$ws.Range("K2").Value = 150737.64
$ws.Range("M2").Value = -150625.64
$ws.Range("H64").Value = 58197.5
$ws.Range("H67").Value = 58197.5
$ws.Range("H116").Value = 68332.336
$ws.Range("J116").Value = 68332.336
$ws.Range("L116").Value = 68332.336
$ws.Range("N116").Value = -77510.336
$ws.Range("H119").ClearContents()
$ws.Range("I119").ClearContents()
$ws.Range("J119").ClearContents()
$ws.Range("K119").ClearContents()
$ws.Range("L119").ClearContents()
$ws.Range("N119").ClearContents()
$ws.Range("H120").ClearContents()
$ws.Range("I120").ClearContents()
$ws.Range("J120").ClearContents()
$ws.Range("K120").ClearContents()
$ws.Range("L120").ClearContents()
$ws.Range("N120").ClearContents()
$ws.Range("H121").ClearContents()
$ws.Range("I121").ClearContents()
$ws.Range("J121").ClearContents()
$ws.Range("K121").ClearContents()
$ws.Range("L121").ClearContents()
$ws.Range("H122").ClearContents()
$ws.Range("I122").ClearContents()
$ws.Range("J122").ClearContents()
$ws.Range("K122").ClearContents()
$ws.Range("L122").ClearContents()
$ws.Range("H123").ClearContents()
$ws.Range("I123").ClearContents()
$ws.Range("J123").ClearContents()
$ws.Range("K123").ClearContents()
$ws.Range("L123").ClearContents()
$ws.Range("H124").ClearContents()
$ws.Range("I124").ClearContents()
$ws.Range("J124").ClearContents()
$ws.Range("K124").ClearContents()
$ws.Range("L124").ClearContents()
$ws.Range("N124").ClearContents()
$ws.Range("H125").ClearContents()
$ws.Range("I125").ClearContents()
$ws.Range("J125").ClearContents()
$ws.Range("K125").ClearContents()
$ws.Range("L125").ClearContents()
$ws.Range("H126").ClearContents()
$ws.Range("I126").ClearContents()
$ws.Range("J126").ClearContents()
$ws.Range("K126").ClearContents()
$ws.Range("L126").ClearContents()
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H127").ClearContents()
$ws.Range("I127").ClearContents()
$ws.Range("J127").ClearContents()
$ws.Range("K127").ClearContents()
$ws.Range("L127").ClearContents()
$ws.Range("H128").ClearContents()
$ws.Range("I128").ClearContents()
$ws.Range("J128").ClearContents()
$ws.Range("K128").ClearContents()
$ws.Range("L128").ClearContents()
$ws.Range("H129").ClearContents()
$ws.Range("I129").ClearContents()
$ws.Range("J129").ClearContents()
$ws.Range("K129").ClearContents()
$ws.Range("L129").ClearContents()
$ws.Range("H130").ClearContents()
$ws.Range("I130").ClearContents()
$ws.Range("J130").ClearContents()
$ws.Range("K130").ClearContents()
$ws.Range("L130").ClearContents()
$ws.Range("N130").ClearContents()
$ws.Range("H131").ClearContents()
$ws.Range("I131").ClearContents()
$ws.Range("J131").ClearContents()
$ws.Range("K131").ClearContents()
$ws.Range("L131").ClearContents()
$ws.Range("H132").ClearContents()
$ws.Range("I132").ClearContents()
$ws.Range("J132").ClearContents()
$ws.Range("K132").ClearContents()
$ws.Range("L132").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H133").ClearContents()
$ws.Range("I133").ClearContents()
$ws.Range("J133").ClearContents()
$ws.Range("K133").ClearContents()
$ws.Range("L133").ClearContents()
$ws.Range("H135").ClearContents()
$ws.Range("I135").ClearContents()
$ws.Range("J135").ClearContents()
$ws.Range("K135").ClearContents()
$ws.Range("L135").ClearContents()
$ws.Range("H136").ClearContents()
$ws.Range("I136").ClearContents()
$ws.Range("J136").ClearContents()
$ws.Range("K136").ClearContents()
$ws.Range("L136").ClearContents()
$ws.Range("M136").ClearContents()
$ws.Range("H137").ClearContents()
$ws.Range("I137").ClearContents()
$ws.Range("J137").ClearContents()
$ws.Range("K137").ClearContents()
$ws.Range("L137").ClearContents()
$ws.Range("N137").ClearContents()
$ws.Range("H138").ClearContents()
$ws.Range("I138").ClearContents()
$ws.Range("J138").ClearContents()
$ws.Range("K138").ClearContents()
$ws.Range("L138").ClearContents()
$ws.Range("H139").ClearContents()
$ws.Range("I139").ClearContents()
$ws.Range("J139").ClearContents()
$ws.Range("K139").ClearContents()
$ws.Range("L139").ClearContents()
$ws.Range("H140").ClearContents()
$ws.Range("I140").ClearContents()
$ws.Range("J140").ClearContents()
$ws.Range("K140").ClearContents()
$ws.Range("L140").ClearContents()
$ws.Range("H141").ClearContents()
$ws.Range("I141").ClearContents()
$ws.Range("J141").ClearContents()
$ws.Range("K141").ClearContents()
$ws.Range("L141").ClearContents()
$ws.Range("N141").ClearContents()

Write-Host "Applied all Sophia_Profits updates"